$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 40, shifting existing rows 40-148 down to 41-149
$ws.Rows("40:40").Insert()

# Populate the newly inserted row 40 with the new weekly record
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = "Vega Modelo de Temuco"
$ws.Range("C40").Value = "La Araucanía"
$ws.Range("D40").Value = 44715
$ws.Range("E40").Value = 9
$ws.Range("F40").Value = 100112012
$ws.Range("G40").Value = "Espinaca"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 65
$ws.Range("K40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = 10000
$ws.Range("N40").Value = "$/docena de atados"
$ws.Range("O40").Value = "Región de La Araucanía"
$ws.Range("P40").Value = 3333
$ws.Range("Q40").Value = 3
$ws.Range("R40").Value = "Hortaliza"
